$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# This document has a single paragraph (containing the "_GoBack" bookmark)
# that reads:
#   "Google Docs originated when Google acquired Upstartle in 2006 and
#   through 2007 merged their web-based word processor with Google
#   Spreadsheets [18][BOOKMARK]. The original intentions, stemming from
#   Upstartle, were to allow users to share documents instantly, as well
#   as collaborate in real-time. This continued with the addition of
#   presentations to Google Docs, and was refined over time to include
#   more features. "
#
# The edit:
#   1. Moves "The original intentions... more features." text to BEFORE
#      the bookmark (appended right after "...Google Spreadsheets [18]"),
#      changing the ending to "...to include extra usability features."
#   2. Splits that text into its own paragraph (bookmark now sits alone,
#      at the start of the next paragraph).
#   3. Adds a blank paragraph, then a new paragraph about Stack Overflow
#      (with "Stack Overflow" in italics), then another blank paragraph,
#      all inserted before the (now solo) bookmark paragraph.
# ---------------------------------------------------------------------------

# --- Step 1: move the "moved" text to before the bookmark, with its new
#     wording, using the bookmark's (zero-width) Range as an anchor so the
#     inserted text lands before the bookmark rather than after it. ---
$bm = $d.Bookmarks("_GoBack")
$r = $bm.Range
$r.InsertBefore(". The original intentions, stemming from Upstartle, were to allow users to share documents instantly, as well as collaborate in real-time. This continued with the addition of presentations to Google Docs, and was refined over time to include extra usability features.")

# --- Step 2: delete the old copy of that text, which is still sitting
#     right after the bookmark. ---
$old = $d.Content
$foundOld = $old.Find.Execute(". The original intentions, stemming from Upstartle, were to allow users to share documents instantly, as well as collaborate in real-time. This continued with the addition of presentations to Google Docs, and was refined over time to include more features. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundOld) {
    $old.Delete()
}

# --- Step 3: split the paragraph right before the bookmark, so the
#     bookmark becomes the sole occupant of its own (new) paragraph. ---
$bm = $d.Bookmarks("_GoBack")
$r = $bm.Range
$r.InsertBefore("`r")

# --- Step 4: insert a blank paragraph before the bookmark paragraph. ---
$bm = $d.Bookmarks("_GoBack")
$r = $bm.Range
$r.InsertBefore("`r")

# --- Step 5: insert the new "Stack Overflow" paragraph before the
#     bookmark paragraph. First insert the italicised lead-in phrase,
#     tracking the exact offsets via the bookmark so we can italicise
#     only that substring afterwards. ---
$bm = $d.Bookmarks("_GoBack")
$startPos = $bm.Start
$r = $bm.Range
$r.InsertBefore("Stack Overflow")
$bm = $d.Bookmarks("_GoBack")
$endPos = $bm.Start
$italicRange = $d.Range($startPos, $endPos)
$italicRange.Font.Italic = $true

# Remainder of the Stack Overflow paragraph (plain formatting).
$bm = $d.Bookmarks("_GoBack")
$r = $bm.Range
$r.InsertBefore(" began in 2008 as a website dedicated to helping users seek assistance on programming related issues. Soon after in 2009, additional websites were created along the same premise, under the Stack Exchange umbrella. Stack Overflow itself continues to play a Q&A facilitator for programmers looking for answers to specific, closed loop questions, and is moderated by the community. ")

# --- Step 6: insert another blank paragraph before the bookmark
#     paragraph. ---
$bm = $d.Bookmarks("_GoBack")
$r = $bm.Range
$r.InsertBefore("`r")

Write-Output "done"
